$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("A3").Value = 131067035
$ws.Range("B3").Value = 79243
$ws.Range("E3").Value = 6425
$ws.Range("F3").Value = "Garnlav"
$ws.Range("G3").Value = "Alectoria sarmentosa"
$ws.Range("H3").Value = "(Ach.) Ach."
$ws.Range("M3").ClearContents()
$ws.Range("Q3").Value = 466172
$ws.Range("R3").Value = 7046340
$ws.Range("AC3").Value = "Långväxta bålar på gran."
$ws.Range("AM3").Value = "Gren på levande träd"
$ws.Range("AO3").Value = "Branch on living tree # Picea abies"

# Row 4
$ws.Range("A4").Value = 131067787
$ws.Range("Q4").Value = 466335
$ws.Range("R4").Value = 7046445
$ws.Range("AH4").Value = "Barrskog"

# Row 5
$ws.Range("A5").Value = 131067798
$ws.Range("Q5").Value = 466279
$ws.Range("R5").Value = 7046403
$ws.Range("AH5").Value = "Granskog"

# Row 6
$ws.Range("A6").Value = 131067788
$ws.Range("B6").Value = 57884
$ws.Range("E6").Value = 100109
$ws.Range("F6").Value = "Tretåig hackspett"
$ws.Range("G6").Value = "Picoides tridactylus"
$ws.Range("H6").Value = "(Linnaeus, 1758)"
$ws.Range("I6").ClearContents()
$ws.Range("M6").Value = "äldre spår"
$ws.Range("N6").ClearContents()
$ws.Range("Q6").Value = 466325
$ws.Range("R6").Value = 7046442
$ws.Range("AC6").Value = "Ringhack, äldre, på gran."
$ws.Range("AJ6").Value = "gran"
$ws.Range("AK6").Value = "Picea abies"
$ws.Range("AO6").Value = "Picea abies"

# Row 7
$ws.Range("A7").Value = 131067810
$ws.Range("B7").Value = 58043
$ws.Range("E7").Value = 103021
$ws.Range("F7").Value = "Talltita"
$ws.Range("G7").Value = "Poecile montanus"
$ws.Range("H7").Value = "(Conrad von Baldenstein, 1827)"
$ws.Range("I7").Value = "2"
$ws.Range("M7").Value = "födosökande"
$ws.Range("N7").Value = "observerad"
$ws.Range("Q7").Value = 466532
$ws.Range("R7").Value = 7046605
$ws.Range("AC7").Value = "Två födosökande talltitor i äldre barrskog med flerskiktning och murknande björkhögstubbar för bohål."
$ws.Range("AJ7").ClearContents()
$ws.Range("AK7").ClearContents()
$ws.Range("AM7").ClearContents()
$ws.Range("AO7").ClearContents()

# Row 10
$ws.Range("A10").Value = 131067030
$ws.Range("B10").Value = 79243
$ws.Range("E10").Value = 6425
$ws.Range("F10").Value = "Garnlav"
$ws.Range("G10").Value = "Alectoria sarmentosa"
$ws.Range("H10").Value = "(Ach.) Ach."
$ws.Range("M10").ClearContents()
$ws.Range("Q10").Value = 466302
$ws.Range("R10").Value = 7046517
$ws.Range("AC10").ClearContents()
$ws.Range("AM10").Value = "Gren på levande träd"
$ws.Range("AO10").Value = "Branch on living tree # Picea abies"

# Row 11
$ws.Range("A11").Value = 131067786
$ws.Range("M11").Value = "äldre spår"
$ws.Range("Q11").Value = 466366
$ws.Range("R11").Value = 7046466
$ws.Range("AC11").Value = "Ringhack, äldre, på gran."
$ws.Range("AM11").ClearContents()
$ws.Range("AO11").Value = "Picea abies"

# Row 12
$ws.Range("A12").Value = 131067792
$ws.Range("M12").Value = "färska spår"
$ws.Range("Q12").Value = 466356
$ws.Range("R12").Value = 7046460
$ws.Range("AC12").Value = "Ringhack, färska, på gran."
$ws.Range("AH12").Value = "Barrskog"
$ws.Range("AM12").Value = "Trädstam på levande träd"
$ws.Range("AO12").Value = "Stem on living tree # Picea abies"

# Row 13
$ws.Range("A13").Value = 131067781
$ws.Range("B13").Value = 57884
$ws.Range("E13").Value = 100109
$ws.Range("F13").Value = "Tretåig hackspett"
$ws.Range("G13").Value = "Picoides tridactylus"
$ws.Range("H13").Value = "(Linnaeus, 1758)"
$ws.Range("M13").Value = "äldre spår"
$ws.Range("Q13").Value = 466204
$ws.Range("R13").Value = 7046448
$ws.Range("AC13").Value = "Ringhack, äldre, på gran."
$ws.Range("AH13").Value = "Granskog"
$ws.Range("AM13").ClearContents()
$ws.Range("AO13").Value = "Picea abies"

# Row 27
$ws.Range("A27").Value = 131067797
$ws.Range("B27").Value = 57884
$ws.Range("E27").Value = 100109
$ws.Range("F27").Value = "Tretåig hackspett"
$ws.Range("G27").Value = "Picoides tridactylus"
$ws.Range("H27").Value = "(Linnaeus, 1758)"
$ws.Range("M27").Value = "äldre spår"
$ws.Range("Q27").Value = 466283
$ws.Range("R27").Value = 7046407
$ws.Range("AC27").Value = "Ringhack, äldre, på gran."
$ws.Range("AM27").ClearContents()
$ws.Range("AO27").Value = "Picea abies"

# Row 28
$ws.Range("A28").Value = 131067038
$ws.Range("B28").Value = 79243
$ws.Range("E28").Value = 6425
$ws.Range("F28").Value = "Garnlav"
$ws.Range("G28").Value = "Alectoria sarmentosa"
$ws.Range("H28").Value = "(Ach.) Ach."
$ws.Range("M28").ClearContents()
$ws.Range("Q28").Value = 466024
$ws.Range("R28").Value = 7046276
$ws.Range("AC28").Value = "Enstaka bålar på gran."
$ws.Range("AM28").Value = "Gren på levande träd"
$ws.Range("AO28").Value = "Branch on living tree # Picea abies"

# Row 37
$ws.Range("A37").Value = 131067790
$ws.Range("B37").Value = 57884
$ws.Range("E37").Value = 100109
$ws.Range("F37").Value = "Tretåig hackspett"
$ws.Range("G37").Value = "Picoides tridactylus"
$ws.Range("H37").Value = "(Linnaeus, 1758)"
$ws.Range("M37").Value = "äldre spår"
$ws.Range("Q37").Value = 466313
$ws.Range("R37").Value = 7046432
$ws.Range("AC37").Value = "Ringhack, äldre, på gran."
$ws.Range("AH37").Value = "Barrskog"
$ws.Range("AM37").ClearContents()
$ws.Range("AO37").Value = "Picea abies"

# Row 38
$ws.Range("A38").Value = 131067789
$ws.Range("B38").Value = 57884
$ws.Range("E38").Value = 100109
$ws.Range("F38").Value = "Tretåig hackspett"
$ws.Range("G38").Value = "Picoides tridactylus"
$ws.Range("H38").Value = "(Linnaeus, 1758)"
$ws.Range("M38").Value = "äldre spår"
$ws.Range("Q38").Value = 466317
$ws.Range("R38").Value = 7046435
$ws.Range("AC38").Value = "Ringhack, äldre, på gran."
$ws.Range("AH38").Value = "Barrskog"
$ws.Range("AM38").ClearContents()
$ws.Range("AO38").Value = "Picea abies"

# Row 39
$ws.Range("A39").Value = 131067031
$ws.Range("B39").Value = 79243
$ws.Range("E39").Value = 6425
$ws.Range("F39").Value = "Garnlav"
$ws.Range("G39").Value = "Alectoria sarmentosa"
$ws.Range("H39").Value = "(Ach.) Ach."
$ws.Range("M39").ClearContents()
$ws.Range("Q39").Value = 466288
$ws.Range("R39").Value = 7046458
$ws.Range("AC39").Value = "På flera gamla granar i granskog."
$ws.Range("AH39").Value = "Granskog"
$ws.Range("AM39").Value = "Gren på levande träd"
$ws.Range("AO39").Value = "Branch on living tree # Picea abies"

# Row 40
$ws.Range("A40").Value = 131067032
$ws.Range("B40").Value = 79243
$ws.Range("E40").Value = 6425
$ws.Range("F40").Value = "Garnlav"
$ws.Range("G40").Value = "Alectoria sarmentosa"
$ws.Range("H40").Value = "(Ach.) Ach."
$ws.Range("M40").ClearContents()
$ws.Range("Q40").Value = 466239
$ws.Range("R40").Value = 7046392
$ws.Range("AC40").ClearContents()
$ws.Range("AH40").Value = "Granskog"
$ws.Range("AM40").Value = "Gren på levande träd"
$ws.Range("AO40").Value = "Branch on living tree # Picea abies"
